# The commit swaps the presentation's main theme color palette from the
# custom "Integral" (Red Violet) scheme over to the stock PowerPoint
# "Office Theme" palette (dk1..folHlink), while font/format schemes stay
# the same. Drive this the way a user would from the Design tab: rewrite
# the 12 theme colors on the (single) slide master's theme color scheme.
#
# Office Theme color values (clrScheme order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), expressed as the BGR-packed decimal .RGB values the
# PowerPoint object model expects:
#   dk1      = 000000 -> 0
#   lt1      = FFFFFF -> 16777215
#   dk2      = 44546A -> 6968388
#   lt2      = E7E6E6 -> 15132391
#   accent1  = 5B9BD5 -> 13998939
#   accent2  = ED7D31 -> 3243501
#   accent3  = A5A5A5 -> 10855845
#   accent4  = FFC000 -> 49407
#   accent5  = 4472C4 -> 12874308
#   accent6  = 70AD47 -> 4697456
#   hlink    = 0563C1 -> 12673797
#   folHlink = 954F72 -> 7491477

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
